$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# DriverSheet: insert a new "GetPassScreenshot" flag column right after
# the existing "Execute" column, mirroring the Y/N values already in A.
# ---------------------------------------------------------------------
$driver = $wb.Worksheets.Item("DriverSheet")
$driver.Columns("B:B").Insert()

$driver.Range("B1").Value = "GetPassScreenshot"
$driver.Range("B2").Value = "Y"
$driver.Range("B3").Value = "Y"
$driver.Range("B4").Value = "Y"
$driver.Range("B5").Value = "Y"
$driver.Range("B6").Value = "N"
$driver.Range("B7").Value = "N"
$driver.Range("B8").Value = "N"
$driver.Range("B9").Value = "N"
$driver.Range("B10").Value = "N"

$driver.Range("B7:B10").Select()
$driver.Activate()

# ---------------------------------------------------------------------
# DataSheet: insert a new "HeadLess" flag column right after the
# existing "execute" column.
# ---------------------------------------------------------------------
$data = $wb.Worksheets.Item("DataSheet")
$data.Columns("B:B").Insert()

$data.Range("B1").Value = "HeadLess"
$data.Range("B2").Value = "F"
$data.Range("B3").Value = "T"
$data.Range("B4").Value = "F"
$data.Range("B5").Value = "F"
$data.Range("B6").Value = "T"
$data.Range("B7").Value = "T"
$data.Range("B8").Value = "T"
$data.Range("B9").Value = "T"

# Existing row toggled from N to Y.
$data.Range("A3").Value = "Y"

$data.Range("D15").Select()
